$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 6275
$ws.Range("I4").Value = 6275
$ws.Range("K4").Value = 6275
$ws.Range("M4").Value = -6161
$ws.Range("H5").Value = 262.83334
$ws.Range("I5").Value = 262.83334
$ws.Range("K5").Value = 262.83334
$ws.Range("M5").Value = -147.83334
$ws.Range("H9").Value = 263.94446
$ws.Range("I9").Value = 178.3125
$ws.Range("J9").Value = 949
$ws.Range("K9").Value = 178.3125
$ws.Range("L9").Value = 949
$ws.Range("M9").Value = -9.3125
$ws.Range("N9").Value = -1287
$ws.Range("H18").Value = 13499
$ws.Range("J18").Value = 1373.75
$ws.Range("L18").Value = 1373.75
$ws.Range("N18").Value = -1941.75
$ws.Range("H19").Value = 1534
$ws.Range("I19").Value = 1598.6666
$ws.Range("J19").Value = 1506.2858
$ws.Range("K19").Value = 1598.6666
$ws.Range("L19").Value = 1506.2858
$ws.Range("M19").Value = -1423.6666
$ws.Range("N19").Value = -1856.2858
$ws.Range("H40").Value = 9375
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 9375
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 9375
$ws.Range("N40").Value = -9725
$ws.Range("H42").Value = 18.75
$ws.Range("I42").Value = 10.333333
$ws.Range("J42").Value = 23.8
$ws.Range("K42").Value = 30.999999
$ws.Range("L42").Value = 71.40000000000001
$ws.Range("M42").Value = 199.000001
$ws.Range("N42").Value = -531.4
$ws.Range("H43").Value = 3399.5
$ws.Range("J43").Value = 3399
$ws.Range("L43").Value = 3399
$ws.Range("N43").Value = -3537
$ws.Range("H69").Value = 2501
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 2501
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 5428.5713
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064
$ws.Range("H77").Value = 5428.5713
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320
$ws.Range("H80").Value = 13123.5
$ws.Range("I80").Value = 495
$ws.Range("J80").Value = 17333
$ws.Range("K80").Value = 1485
$ws.Range("L80").Value = 51999
$ws.Range("M80").Value = -487
$ws.Range("N80").Value = -53995
$ws.Range("H83").Value = 13123.5
$ws.Range("I83").Value = 495
$ws.Range("J83").Value = 17333
$ws.Range("K83").Value = 4455
$ws.Range("L83").Value = 155997
$ws.Range("M83").Value = 537
$ws.Range("N83").Value = -165981
$ws.Range("H100").Value = 6626.25
$ws.Range("I100").Value = 6502.5
$ws.Range("J100").Value = 6750
$ws.Range("K100").Value = 6502.5
$ws.Range("L100").Value = 6750
$ws.Range("M100").Value = -5961.5
$ws.Range("N100").Value = -7832
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H130").Value = 94777.8
$ws.Range("J130").Value = 94777.8
$ws.Range("L130").Value = 94777.8
$ws.Range("N130").Value = -104817.8
$ws.Range("H132").Value = 3786.5557
$ws.Range("I132").Value = 3402.3333
$ws.Range("K132").Value = 10206.9999
$ws.Range("M132").Value = -7676.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2549.75
$ws.Range("J97").Value = 7800
$ws.Range("L97").Value = 7800
$ws.Range("N97").Value = -8792
$ws.Range("H101").Value = 12997.5
$ws.Range("J101").Value = 12997.5
$ws.Range("L101").Value = 12997.5
$ws.Range("N101").Value = -19487.5
$ws.Range("H132").Value = 1729.3334
$ws.Range("I132").Value = 1729.3334
$ws.Range("K132").Value = 5188.0002
$ws.Range("M132").Value = -2658.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1999.5
$ws.Range("I107").Value = 1999.5
$ws.Range("K107").Value = 1999.5
$ws.Range("M107").Value = -79.5
$ws.Range("H134").Value = 6560.636
$ws.Range("I134").Value = 4234.2
$ws.Range("K134").Value = 12702.6
$ws.Range("M134").Value = -10167.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 92.5
$ws.Range("I7").Value = 76.85714
$ws.Range("J7").Value = 129
$ws.Range("K7").Value = 76.85714
$ws.Range("L7").Value = 129
$ws.Range("M7").Value = 36.14286
$ws.Range("N7").Value = -355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 400601.4
$ws.Range("I4").Value = 417251.9
$ws.Range("J4").Value = 333999.34
$ws.Range("K4").Value = 1251755.7
$ws.Range("L4").Value = 1001998.02
$ws.Range("M4").Value = -1251643.7
$ws.Range("N4").Value = -1002222.02
$ws.Range("H17").Value = 203.33333
$ws.Range("I17").Value = 157.5
$ws.Range("J17").Value = 295
$ws.Range("K17").Value = 472.5
$ws.Range("L17").Value = 885
$ws.Range("M17").Value = -303.5
$ws.Range("N17").Value = -1223

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3209.1
$ws.Range("I80").Value = 2756.2856
$ws.Range("K80").Value = 2756.2856
$ws.Range("M80").Value = -1758.2856
$ws.Range("H83").Value = 3209.1
$ws.Range("I83").Value = 2756.2856
$ws.Range("K83").Value = 13781.428
$ws.Range("M83").Value = -8789.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 3824.182
$ws.Range("I40").Value = 3563
$ws.Range("K40").Value = 3563
$ws.Range("M40").Value = -3427
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("N41").Value = 0
$ws.Range("H46").Value = 3632.75
$ws.Range("J46").Value = 4468.6665
$ws.Range("L46").Value = 4468.6665
$ws.Range("N46").Value = -4844.6665
$ws.Range("H100").Value = 4554.6665
$ws.Range("I100").Value = 4554.6665
$ws.Range("K100").Value = 4554.6665
$ws.Range("M100").Value = -4013.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("M54").Value = 20000
$ws.Range("N54").Value = -21040
$ws.Range("H81").Value = 17588.445
$ws.Range("I81").Value = 17588.445
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 35176.89
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -34115.89
$ws.Range("H84").Value = 17588.445
$ws.Range("I84").Value = 17588.445
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 175884.45
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -170580.45
$ws.Range("H107").Value = 1864.1818
$ws.Range("I107").Value = 1255.8889
$ws.Range("K107").Value = 3767.6667
$ws.Range("M107").Value = -1847.6667
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("H126").Value = 1298.8
$ws.Range("J126").Value = 1499.6666
$ws.Range("L126").Value = 4498.9998
$ws.Range("N126").Value = -9438.9998
